# Build site at 2023-04-12 14:53:07 UTC
# Apply the LOT2007 content restructuring: new Objetivos text, new
# Docentes/Programa resumido/Programa blocks inserted, Avaliacao block
# shifted down, new Bibliografia text, and an extra Requisitos row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Objetivos: -> new objectives paragraph (PT) --------------
$ws.Range("B10").Value = 'Promover aos participantes do curso conhecimentos de bioquímica  abrangendo a organização estrutural e molecular da célulaCompreender a importância dos compostos orgânicos no metabolismo celularUtilizar todos os conhecimentos como pré-requisito para as disciplinas do curso de engenharia Bioquímica'
$ws.Range("C10").Value = 'Promover aos participantes do curso conhecimentos de bioquímica  abrangendo a organização estrutural e molecular da célulaCompreender a importância dos compostos orgânicos no metabolismo celularUtilizar todos os conhecimentos como pré-requisito para as disciplinas do curso de engenharia Bioquímica'

# --- Row 11 (Objectives:) and Row 12 (Docentes responsáveis:) stay ----
# (unchanged, nothing to do)

# --- Row 13: used to be "Programa resumido: / Semestral", now holds --
# the professor name that used to live in row 18, with no label in A
# and the default (non custom) row height.
$ws.Range("A13").Clear()
$ws.Range("B13").Value = '427823 - Adriane Maria Ferreira Milagres'
$ws.Range("C13").Value = '427823 - Adriane Maria Ferreira Milagres'
$ws.Rows.Item(13).AutoFit()

# --- Row 14: now "Programa resumido:" with new summarized syllabus ---
$ws.Range("A14").Value = 'Programa resumido:'
$ws.Range("B14").Value = '01Química ácido-base/Tampões02Aminoácidos03Proteínas:Estrutura primária04Proteínas:Estrutura tridimensional05Função das proteínas06Enzimas: catálise enzimática07Cinética enzimática, inibição e regulação08Carboidratos09Lipídeos10Membranas Biológicas11Nucleotídeos e ácidos nucleicos'
$ws.Range("C14").Value = '01Química ácido-base/Tampões02Aminoácidos03Proteínas:Estrutura primária04Proteínas:Estrutura tridimensional05Função das proteínas06Enzimas: catálise enzimática07Cinética enzimática, inibição e regulação08Carboidratos09Lipídeos10Membranas Biológicas11Nucleotídeos e ácidos nucleicos'
$ws.Rows.Item(14).RowHeight = 60

# --- Row 15: now "Short syllabus:" with the English short syllabus ---
$ws.Range("A15").Value = 'Short syllabus:'
$ws.Range("B15").Value = '1.Acids and bases/Buffer solutions; 2. Amino acids; 3. Proteins: The primary level of protein structure; 4. Proteins: The three-dimensional structure; 5. Protein Function; 6. Enzymes: Biological catalysts; 7. The kinetics of enzymatic catalysis; 8. Carbohydrates; 9. Lipids; 10. Membranes and cellular transport; 11. Nucleic acids.'
$ws.Range("C15").Value = '1.Acids and bases/Buffer solutions; 2. Amino acids; 3. Proteins: The primary level of protein structure; 4. Proteins: The three-dimensional structure; 5. Protein Function; 6. Enzymes: Biological catalysts; 7. The kinetics of enzymatic catalysis; 8. Carbohydrates; 9. Lipids; 10. Membranes and cellular transport; 11. Nucleic acids.'
$ws.Rows.Item(15).RowHeight = 60

# --- Row 16: now "Programa:" with the new long PT syllabus -----------
$ws.Range("A16").Value = 'Programa:'
$ws.Range("B16").Value = '01Química ácido-base/Tampões : Constante de dissociação, Curvas de titulação , capacidade tamponante02Aminoácidos:Estrutura dos aminoácidos, classificação e características, Nomenclatura, propriedades ácido-base, estereoquímica,aminoácidos incomuns03Proteínas:Estrutura primáriaPurificação de proteínas, solubilidade, cormatografia, eletroforese. Sequenciamento de proteínas.04Proteínas:Estrutura tridimensionalEstrutura secundária, terciária, quaternária. Dobramento e estabilidade das proteínas.05Função das proteínasMioglobina, hemoglobina, anticorpos06Enzimas: catálise enzimáticaNomenclatura das enzimas, especificidade dos substratos, co-fatores e coenzimas, Energia de ativação e coordenada de reação. Curvas de progresso. Efeito da temperatura e pH sobre a velocidade das reações enzimáticas07Cinética enzimática, inibição e regulação Efeito da concentração de substrato na velocidade das reações enzimáticas. Efeito de inibidores. Inibições reversíveis.  Modelos de inibição competitiva, não competitiva e acompetitiva simples.08CarboidratosMonossacarídeos:classificação, configuração e conformação. Dissacarídeos, Polissacarídeos estruturais:celulose e quitina, Polissacarídeos de reserva:amido e glicogênio, glicosaminoglicanos, Glicoproteínas:oligossacarídeos, paredes celulares bacterianas.09LipídeosClassificação: ácidos graxos, trioacilglicerol, glicerofosfolipídeos, esfingolipídeos, esteróides.10Membranas BiológicasProteínas de membrana:integrais e periféricas, modelo do mosaico fluido, assimetria dos lipídeos, Transporte através da membrana: termodinâmica do transporte,  transporte passivo e ativo.11Nucleotídeos e ácidos nucleicos:Estrutura e função dos nucleotídeos, Estrutura dos ácidos nucleicos, sequenciamento de ácidos nucleicos, endonucleases de restrição, Bibliotecas genômicas,Amplificação do DNA pela reação em cadeia da polimerase.'
$ws.Range("C16").Value = '01Química ácido-base/Tampões : Constante de dissociação, Curvas de titulação , capacidade tamponante02Aminoácidos:Estrutura dos aminoácidos, classificação e características, Nomenclatura, propriedades ácido-base, estereoquímica,aminoácidos incomuns03Proteínas:Estrutura primáriaPurificação de proteínas, solubilidade, cormatografia, eletroforese. Sequenciamento de proteínas.04Proteínas:Estrutura tridimensionalEstrutura secundária, terciária, quaternária. Dobramento e estabilidade das proteínas.05Função das proteínasMioglobina, hemoglobina, anticorpos06Enzimas: catálise enzimáticaNomenclatura das enzimas, especificidade dos substratos, co-fatores e coenzimas, Energia de ativação e coordenada de reação. Curvas de progresso. Efeito da temperatura e pH sobre a velocidade das reações enzimáticas07Cinética enzimática, inibição e regulação Efeito da concentração de substrato na velocidade das reações enzimáticas. Efeito de inibidores. Inibições reversíveis.  Modelos de inibição competitiva, não competitiva e acompetitiva simples.08CarboidratosMonossacarídeos:classificação, configuração e conformação. Dissacarídeos, Polissacarídeos estruturais:celulose e quitina, Polissacarídeos de reserva:amido e glicogênio, glicosaminoglicanos, Glicoproteínas:oligossacarídeos, paredes celulares bacterianas.09LipídeosClassificação: ácidos graxos, trioacilglicerol, glicerofosfolipídeos, esfingolipídeos, esteróides.10Membranas BiológicasProteínas de membrana:integrais e periféricas, modelo do mosaico fluido, assimetria dos lipídeos, Transporte através da membrana: termodinâmica do transporte,  transporte passivo e ativo.11Nucleotídeos e ácidos nucleicos:Estrutura e função dos nucleotídeos, Estrutura dos ácidos nucleicos, sequenciamento de ácidos nucleicos, endonucleases de restrição, Bibliotecas genômicas,Amplificação do DNA pela reação em cadeia da polimerase.'
$ws.Rows.Item(16).RowHeight = 120

# --- Row 17: now "Syllabus:" with the English long syllabus ----------
$ws.Range("A17").Value = 'Syllabus:'
$ws.Range("B17").Value = '1.Acid-bases chemistry/Buffers: dissociation constant, titration curves, buffering capacity. Ionisation equilibria of acids and bases in aqueous solutions2.Amino acids: structure of the amino acids, properties of amino acids side chains, classes and nomenclature, acid-bases properties, stereochemistry, Modified amino acids.3.Proteins: primary structure, peptides and the peptide bond. Protein purification, solubility, chromatography, electrophoresis. Proteins sequences.4.Proteins: three-dimensional structure. Secondary structure, tertiary structure, quaternary structure. Dynamics of protein structure:  folding and stability.5.Protein Function: oxygen-binding proteins: myoglobin and hemoglobin, immunoglobulins.6.Enzymes: Biological catalysis. Nomenclature and classification of enzymes, specificity of the substrates, co-factors and co-enzymes. Energy of activation and coordinator of reaction. Progress curves. Effect of temperature and pH on the rate of enzymatic activity. 7.Enzymatic kinetics, inhibition and regulation: Effect of the substrate concentration on the rate of enzymatic reactions. Enzyme Inhibition. Reversible inhibition. Models of competitive, non-competitive and simple competitive inhibition. Allosteric regulation.8.Carbohydrates: Monosaccharides, stereoisomerism, classification, configuration and conformation. Derivatives of monosaccharides, oligosaccharides, structural polysaccharides: cellulose and chitin, storage polysaccharides: starch and glycogen, Glycoconjugates: Proteoglycans, Glycoproteins, and Glycolipids.9.Lipids: Classification, fatty acids, triacylglycerols, waxes, glycerophospholipids, esphingolipids, cholesterol.10.Biological membranes. Proteins in membranes: integral and peripheral, the fluid mosaic models, the asymmetry of membranes, transport across membranes: the thermodynamics of transport, passive and active transport.11.Nucleotides and nucleic acids: structure and function of nucleotides. Primary structure of nucleic acids, nucleic acid sequencing, restriction endonucleases. Genome sequency. DNA amplification by the polymerase chain reaction.'
$ws.Range("C17").Value = '1.Acid-bases chemistry/Buffers: dissociation constant, titration curves, buffering capacity. Ionisation equilibria of acids and bases in aqueous solutions2.Amino acids: structure of the amino acids, properties of amino acids side chains, classes and nomenclature, acid-bases properties, stereochemistry, Modified amino acids.3.Proteins: primary structure, peptides and the peptide bond. Protein purification, solubility, chromatography, electrophoresis. Proteins sequences.4.Proteins: three-dimensional structure. Secondary structure, tertiary structure, quaternary structure. Dynamics of protein structure:  folding and stability.5.Protein Function: oxygen-binding proteins: myoglobin and hemoglobin, immunoglobulins.6.Enzymes: Biological catalysis. Nomenclature and classification of enzymes, specificity of the substrates, co-factors and co-enzymes. Energy of activation and coordinator of reaction. Progress curves. Effect of temperature and pH on the rate of enzymatic activity. 7.Enzymatic kinetics, inhibition and regulation: Effect of the substrate concentration on the rate of enzymatic reactions. Enzyme Inhibition. Reversible inhibition. Models of competitive, non-competitive and simple competitive inhibition. Allosteric regulation.8.Carbohydrates: Monosaccharides, stereoisomerism, classification, configuration and conformation. Derivatives of monosaccharides, oligosaccharides, structural polysaccharides: cellulose and chitin, storage polysaccharides: starch and glycogen, Glycoconjugates: Proteoglycans, Glycoproteins, and Glycolipids.9.Lipids: Classification, fatty acids, triacylglycerols, waxes, glycerophospholipids, esphingolipids, cholesterol.10.Biological membranes. Proteins in membranes: integral and peripheral, the fluid mosaic models, the asymmetry of membranes, transport across membranes: the thermodynamics of transport, passive and active transport.11.Nucleotides and nucleic acids: structure and function of nucleotides. Primary structure of nucleic acids, nucleic acid sequencing, restriction endonucleases. Genome sequency. DNA amplification by the polymerase chain reaction.'
$ws.Rows.Item(17).RowHeight = 120

# --- Row 18: now just "Avaliação:" (the method text moves to row 19) -
$ws.Range("A18").Value = 'Avaliação:'
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()
$ws.Rows.Item(18).AutoFit()

# --- Row 19: now "Método:" / avaliação escrita text -------------------
$ws.Range("A19").Value = 'Método:'
$ws.Range("B19").Value = 'A avaliação será feita por meio de provas escritas.'
$ws.Range("C19").Value = 'A avaliação será feita por meio de provas escritas.'
$ws.Rows.Item(19).RowHeight = 60

# --- Row 20: now "Critério:" / nota final formula ---------------------
$ws.Range("A20").Value = 'Critério:'
$ws.Range("B20").Value = 'A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + 2*P2)/3'
$ws.Range("C20").Value = 'A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + 2*P2)/3'
$ws.Rows.Item(20).RowHeight = 60

# --- Row 21: now "Norma de recuperação:" / recuperação text -----------
$ws.Range("A21").Value = 'Norma de recuperação:'
$ws.Range("B21").Value = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2'
$ws.Range("C21").Value = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2'
$ws.Rows.Item(21).RowHeight = 60

# --- Row 22: now "Bibliografia:" with the two-reference list ----------
$bibliografia = "1. M. Cox, Michael; Nelson, David L.Princípios de Bioquímica de Lehninger - Editora Artmed 6ª Ed. 2014`n2. Voet, D., Voet, J. G., Pratt, C.W. Fundamentos de Bioquímica:a vida em nivel molecular  Editora Artmed, 2014"
$ws.Range("A22").Value = 'Bibliografia:'
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia
$ws.Rows.Item(22).RowHeight = 120

# --- Row 23: now just "Requisitos:" (the prerequisite rows move down) -
$ws.Range("A23").Value = 'Requisitos:'
$ws.Range("B23").Clear()
$ws.Range("C23").Clear()
$ws.Rows.Item(23).AutoFit()

# --- Row 24: first prerequisite (was row 23) ---------------------------
$ws.Range("A24").Clear()
$ws.Range("B24").Value = 'LOT2002 -  Biologia Celular  (Requisito fraco)' + "`n"
$ws.Range("C24").Value = 'LOT2002 -  Biologia Celular  (Requisito fraco)' + "`n"
$ws.Rows.Item(24).RowHeight = 30

# --- Row 25: second prerequisite (was row 24, now a new extra row) -----
$ws.Range("B25").Value = 'LOT2059 -  Química Orgânica Fundamental  (Requisito fraco)' + "`n"
$ws.Range("C25").Value = 'LOT2059 -  Química Orgânica Fundamental  (Requisito fraco)' + "`n"
$ws.Rows.Item(25).RowHeight = 30
